$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column (D) previously held a date+time serial (6/26/2011) for every
# row, formatted as m/d/yy h:mm (numFmtId 22). Replace it with a distinct,
# date-only value per row (1/1/2000 .. 1/5/2000) formatted as a plain date
# (numFmtId 14), and narrow the column to fit the shorter text.

$ws.Range("D2").Value2 = 36526   # 1/1/2000
$ws.Range("D3").Value2 = 36527   # 1/2/2000
$ws.Range("D4").Value2 = 36528   # 1/3/2000
$ws.Range("D5").Value2 = 36529   # 1/4/2000
$ws.Range("D6").Value2 = 36530   # 1/5/2000

$ws.Range("D2:D6").NumberFormat = "mm-dd-yy"

$ws.Columns.Item(4).ColumnWidth = 9.5
